$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows (B:G), rows 2-9
$data = @{
    2  = @(0.2312766622086449, 1.879048156589436, 10.90505104412208, 3.3022796738196, 3.326949562682384, 51)
    3  = @(0.515602234217057, 1.963007146147895, 11.82602839579948, 3.438899300037655, 3.434545796526614, 50)
    4  = @(0.3567007860761084, 1.897578673147804, 11.03219032630682, 3.32147411946967, 3.336486390714962, 49)
    5  = @(0.5452471532615414, 1.991328350292171, 12.46298015738944, 3.530294627561478, 3.524844688185481, 48)
    6  = @(0.4748333535619241, 1.853041692924432, 11.1369772540728, 3.337210999333545, 3.338969495412706, 47)
    7  = @(0.4980755295490984, 1.841574822490943, 11.41347372914766, 3.378383301099457, 3.378389316479261, 46)
    8  = @(0.3675858156243708, 1.629918664393685, 9.8136786659588, 3.132679151454678, 3.14619239584994, 45)
    9  = @(0.3916378424397349, 1.741856079646775, 10.36499202449527, 3.219470767765297, 3.232505401330167, 44)
    10 = @(0.5067708388377236, 1.530135071712105, 9.825491045469901, 3.134563932267119, 3.129936092387531, 43)
    11 = @(0.2841523553027562, 1.690716294712478, 10.17884259426374, 3.190429844748782, 3.216270379803536, 42)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}

# New row 11 label "Q9" in column A, copying style from A10
$ws.Cells.Item(11, 1).Value = "Q9"
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
